$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) stays the same values, columns E/F are Prerequisites/CoRequisites.
# Data rows 2-16: resorted by Course_Code, plus four new Spring 2026 (1261) DMED courses
# appended at the end (DMED560-DMED563).

$rows = @(
    @('DMED500','Course','DMED','DMED 500 - Foundations of Digital Media','None','None','CAT','MDM PROG',''),
    @('DMED501','Course','DMED','DMED 501 - Visual Storytelling','None','None','CAT','MDM PROG',''),
    @('DMED502','Course','DMED','DMED 502 - Improvisation for Collaboration','None','None','CAT','MDM PROG',''),
    @('DMED503','Course','DMED','DMED 503 - Game Design and Gamification','None','None','CAT','MDM PROG',''),
    @('DMED520','Course','DMED','DMED 520 - Projects I: Building Digital Artifacts','None','None','CAT','MDM PROG',''),
    @('DMED521','Course','DMED','DMED 521 - Projects II','DMED520','None','CAT','MDM PROG','REQ-DMED 520 with a grade of B or higher.'),
    @('DMED522','Course','DMED','DMED 522 - Projects III','DMED521','None','CAT','MDM PROG','REQ-DMED 521 with a grade of B or higher.'),
    @('DMED530','Course','DMED','DMED 530 - Internship','None','None','CAT','MDM PROG',''),
    @('DMED531','Course','DMED','DMED 531 - Internship Continuation','DMED530','None','CAT','MDM PROG','REQ: DMED 530.'),
    @('DMED540','Course','DMED','DMED 540 - Special Topics in Digital Media','None','None','CAT','MDM PROG',''),
    @('DMED550','Course','DMED','DMED 550 - Directed Studies in Digital Media','None','None','CAT','MDM PROG',''),
    @('DMED560','Course','DMED','DMED 560 - Production Team Dynamics','None','None','CAT','MDM PROG',''),
    @('DMED561','Course','DMED','DMED 561 - Pitching for Digital Media Professionals I','None','None','CAT','MDM PROG',''),
    @('DMED562','Course','DMED','DMED 562 - Pitching for Digital Media Professionals II','None','None','CAT','MDM PROG',''),
    @('DMED563','Course','DMED','DMED 563 - Multi-Platform Media','None','None','CAT','MDM PROG','')
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r = $r + 1
}

# Column A was widened to fit the longer course codes / header text.
$ws.Columns.Item(1).ColumnWidth = 12.436197916666666

# Selection moved to H9 (matches the saved view in the updated file).
$ws.Range("H9").Select()
